$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.375
$ws.Range("D2").Value = 0.3333333333333333

# Row 3
$ws.Range("B3").Value = 0.5625
$ws.Range("C3").Value = 0.6428571428571429
$ws.Range("D3").Value = 0.6000000000000001

# Row 4
$ws.Range("B4").Value = 0.5
$ws.Range("C4").Value = 0.5
$ws.Range("D4").Value = 0.5
$ws.Range("E4").Value = 0.5

# Row 5
$ws.Range("B5").Value = 0.46875
$ws.Range("C5").Value = 0.4714285714285714
$ws.Range("D5").Value = 0.4666666666666667

# Row 6
$ws.Range("B6").Value = 0.484375
$ws.Range("C6").Value = 0.5
$ws.Range("D6").Value = 0.4888888888888889

# Row 17
$ws.Range("B17").Value = 0.2222222222222222
$ws.Range("C17").Value = 0.2
$ws.Range("D17").Value = 0.2105263157894737

# Row 18
$ws.Range("B18").Value = 0.4666666666666667
$ws.Range("C18").Value = 0.5
$ws.Range("D18").Value = 0.4827586206896552

# Row 19
$ws.Range("B19").Value = 0.375
$ws.Range("C19").Value = 0.375
$ws.Range("D19").Value = 0.375
$ws.Range("E19").Value = 0.375

# Row 20
$ws.Range("B20").Value = 0.3444444444444444
$ws.Range("C20").Value = 0.35
$ws.Range("D20").Value = 0.3466424682395645

# Row 21
$ws.Range("B21").Value = 0.3648148148148149
$ws.Range("C21").Value = 0.375
$ws.Range("D21").Value = 0.3693284936479129

# Row 22
$ws.Range("B22").Value = 0.4545454545454545
$ws.Range("D22").Value = 0.4761904761904762

# Row 23
$ws.Range("B23").Value = 0.6153846153846154
$ws.Range("C23").Value = 0.5714285714285714
$ws.Range("D23").Value = 0.5925925925925927

# Row 24
$ws.Range("B24").Value = 0.5416666666666666
$ws.Range("C24").Value = 0.5416666666666666
$ws.Range("D24").Value = 0.5416666666666666
$ws.Range("E24").Value = 0.5416666666666666

# Row 25
$ws.Range("B25").Value = 0.534965034965035
$ws.Range("C25").Value = 0.5357142857142857
$ws.Range("D25").Value = 0.5343915343915344

# Row 26
$ws.Range("B26").Value = 0.5483682983682984
$ws.Range("C26").Value = 0.5416666666666666
$ws.Range("D26").Value = 0.5440917107583775
